$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.776.81"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "1.855.39"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.031"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07423"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8877"
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.866.80"
$ws.Range("E12").Value = "  +0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.527"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.039"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009094"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.031"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").Value = "27.814.65"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").Value = "2.098.12"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").Value = "  +6.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.005"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.365"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09099"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7749"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.213"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.019"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.033"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.156"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05311"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.864"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5207"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.923"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1675"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.771"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.23%  "

$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.034"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("E48").Value = "  +2.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.715"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.38%  "

$ws.Range("E50").Value = "  +1.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.889"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
